# Apply the "mostly ig generated files" regeneration edit to
# StructureDefinition-end-age.xlsx
#
# Summary of the change:
#  - Metadata sheet: URL value now points at the "cicada" IG instead of "pythia"
#  - Metadata sheet: Date value was regenerated (new timestamp)
#  - Metadata sheet: a new "Jurisdiction" property row was inserted (with an
#    empty value) right after the "Contact" row, pushing Description/Purpose/
#    Copyright/etc. down by one row
#  - Elements sheet: column widths were recalculated (best-fit) because the
#    regenerated content changed

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# 1. Update the URL value (row 2)
$meta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/end-age"

# 2. Update the Date value (row 8)
$meta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# 3. Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" (row 11)
$meta.Rows.Item(11).Insert()

# Copy the formatting of the (now shifted) Description row onto the new row
# so the new cells keep the same style used throughout the table.
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# 4. The Elements sheet content shifts its shared-string references as a
# side effect of the metadata changes above; re-fit its columns since the
# underlying text driving the "best fit" column widths changed.
$elements = $wb.Worksheets.Item("Elements")
$elements.Columns.AutoFit()
